{"js": "// Update answer cells in the division-practice table.\n// Each cell contains a single run with the full \"A\u00f7B=C, D\" string,\n// so a whole-string search-and-replace per cell is sufficient and safe\n// (all old values below are unique within the document).\nconst replacements = [\n  [\"266\u00f73=88, 2\", \"921\u00f79=102, 3\"],\n  [\"220\u00f75=44, 0\", \"212\u00f72=106, 0\"],\n  [\"483\u00f74=120, 3\", \"424\u00f76=70, 4\"],\n  [\"940\u00f79=104, 4\", \"695\u00f77=99, 2\"],\n  [\"923\u00f74=230, 3\", \"793\u00f79=88, 1\"],\n  [\"560\u00f77=80, 0\", \"621\u00f78=77, 5\"],\n  [\"256\u00f76=42, 4\", \"399\u00f72=199, 1\"],\n  [\"158\u00f72=79, 0\", \"787\u00f74=196, 3\"],\n  [\"847\u00f73=282, 1\", \"433\u00f76=72, 1\"],\n  [\"844\u00f76=140, 4\", \"971\u00f76=161, 5\"],\n  [\"936\u00f76=156, 0\", \"961\u00f72=480, 1\"],\n  [\"695\u00f79=77, 2\", \"291\u00f73=97, 0\"],\n  [\"670\u00f72=335, 0\", \"290\u00f76=48, 2\"],\n  [\"536\u00f73=178, 2\", \"568\u00f72=284, 0\"],\n  [\"461\u00f76=76, 5\", \"265\u00f72=132, 1\"],\n  [\"414\u00f76=69, 0\", \"963\u00f73=321, 0\"],\n  [\"101\u00f75=20, 1\", \"136\u00f74=34, 0\"],\n  [\"533\u00f79=59, 2\", \"684\u00f78=85, 4\"],\n  [\"750\u00f73=250, 0\", \"284\u00f73=94, 2\"],\n  [\"765\u00f76=127, 3\", \"155\u00f78=19, 3\"],\n  [\"626\u00f77=89, 3\", \"978\u00f77=139, 5\"],\n  [\"519\u00f73=173, 0\", \"175\u00f72=87, 1\"],\n  [\"701\u00f73=233, 2\", \"456\u00f73=152, 0\"],\n  [\"687\u00f75=137, 2\", \"427\u00f72=213, 1\"],\n  [\"464\u00f76=77, 2\", \"505\u00f77=72, 1\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  // eslint-disable-next-line no-await-in-loop\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update answer cells in the division-practice table.\n# Each cell holds a single run containing the full \"A\u00f7B=C, D\" string,\n# so a whole-string Find/Replace per cell is sufficient and safe\n# (every old value below is unique within the document).\n$d = $word.ActiveDocument\n\n$replacements = @(\n  ,@(\"266\u00f73=88, 2\", \"921\u00f79=102, 3\")\n  ,@(\"220\u00f75=44, 0\", \"212\u00f72=106, 0\")\n  ,@(\"483\u00f74=120, 3\", \"424\u00f76=70, 4\")\n  ,@(\"940\u00f79=104, 4\", \"695\u00f77=99, 2\")\n  ,@(\"923\u00f74=230, 3\", \"793\u00f79=88, 1\")\n  ,@(\"560\u00f77=80, 0\", \"621\u00f78=77, 5\")\n  ,@(\"256\u00f76=42, 4\", \"399\u00f72=199, 1\")\n  ,@(\"158\u00f72=79, 0\", \"787\u00f74=196, 3\")\n  ,@(\"847\u00f73=282, 1\", \"433\u00f76=72, 1\")\n  ,@(\"844\u00f76=140, 4\", \"971\u00f76=161, 5\")\n  ,@(\"936\u00f76=156, 0\", \"961\u00f72=480, 1\")\n  ,@(\"695\u00f79=77, 2\", \"291\u00f73=97, 0\")\n  ,@(\"670\u00f72=335, 0\", \"290\u00f76=48, 2\")\n  ,@(\"536\u00f73=178, 2\", \"568\u00f72=284, 0\")\n  ,@(\"461\u00f76=76, 5\", \"265\u00f72=132, 1\")\n  ,@(\"414\u00f76=69, 0\", \"963\u00f73=321, 0\")\n  ,@(\"101\u00f75=20, 1\", \"136\u00f74=34, 0\")\n  ,@(\"533\u00f79=59, 2\", \"684\u00f78=85, 4\")\n  ,@(\"750\u00f73=250, 0\", \"284\u00f73=94, 2\")\n  ,@(\"765\u00f76=127, 3\", \"155\u00f78=19, 3\")\n  ,@(\"626\u00f77=89, 3\", \"978\u00f77=139, 5\")\n  ,@(\"519\u00f73=173, 0\", \"175\u00f72=87, 1\")\n  ,@(\"701\u00f73=233, 2\", \"456\u00f73=152, 0\")\n  ,@(\"687\u00f75=137, 2\", \"427\u00f72=213, 1\")\n  ,@(\"464\u00f76=77, 2\", \"505\u00f77=72, 1\")\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Execute(\n    $oldText,   # FindText\n    $true,      # MatchCase\n    $false,     # MatchWholeWord\n    $false,     # MatchWildcards\n    $false,     # MatchSoundsLike\n    $false,     # MatchAllWordForms\n    $true,      # Forward\n    1,          # Wrap: wdFindContinue\n    $false,     # Format\n    $newText,   # ReplaceWith\n    2           # Replace: wdReplaceAll\n  ) | Out-Null\n}\n"}
